$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 260 for the newest weekly price update (week of 45127),
# pushing the existing historical rows (260-290) down to (263-293).
$ws.Rows("260:262").Insert()

# New data for the inserted rows, following the same record layout as every
# other row in this "Femacal de La Calera - Chirimoya" sheet.
$newRows = @(
    @{ Row = 260; L = "Especial"; M = 54; N = 32000; S = 3200; R = "Provincia del Elquí" },
    @{ Row = 261; L = "Primera";  M = 48; N = 30000; S = 3000; R = "Provincia del Elquí" },
    @{ Row = 262; L = "Segunda";  M = 35; N = 27000; S = 2700; R = "Provincia del Elquí" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = 45127
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107002
    $ws.Cells.Item($row, 10).Value = "Chirimoya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.N
    $ws.Cells.Item($row, 16).Value = $r.N
    $ws.Cells.Item($row, 17).Value = "$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 10
}
